# Update the date line.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-03-18 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-19 Tuesday", 2)

# Update the division problems in the table. Addressing cells directly by
# (row, column) avoids any ambiguity from duplicate/overlapping text values
# that appear both before and after the edit.
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "51÷6=" },
    @{ Row = 1;  Col = 2; Text = "93÷8=" },
    @{ Row = 1;  Col = 3; Text = "39÷2=" },
    @{ Row = 1;  Col = 4; Text = "20÷6=" },
    @{ Row = 1;  Col = 5; Text = "33÷4=" },

    @{ Row = 5;  Col = 1; Text = "41÷5=" },
    @{ Row = 5;  Col = 2; Text = "26÷8=" },
    @{ Row = 5;  Col = 3; Text = "31÷6=" },
    @{ Row = 5;  Col = 4; Text = "42÷7=" },
    @{ Row = 5;  Col = 5; Text = "64÷5=" },

    @{ Row = 9;  Col = 1; Text = "95÷4=" },
    @{ Row = 9;  Col = 2; Text = "99÷2=" },
    @{ Row = 9;  Col = 3; Text = "29÷3=" },
    @{ Row = 9;  Col = 4; Text = "28÷4=" },
    @{ Row = 9;  Col = 5; Text = "52÷4=" },

    @{ Row = 13; Col = 1; Text = "30÷5=" },
    @{ Row = 13; Col = 2; Text = "60÷8=" },
    @{ Row = 13; Col = 3; Text = "39÷6=" },
    @{ Row = 13; Col = 4; Text = "54÷5=" },
    @{ Row = 13; Col = 5; Text = "62÷3=" },

    @{ Row = 17; Col = 1; Text = "68÷9=" },
    @{ Row = 17; Col = 2; Text = "33÷6=" },
    @{ Row = 17; Col = 3; Text = "20÷2=" },
    @{ Row = 17; Col = 4; Text = "55÷9=" },
    @{ Row = 17; Col = 5; Text = "34÷5=" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
